# Fix bibliotheca sql write error: the per-currency report rows were written
# in an unsorted (SQL result) order; rewrite rows 2-23 sorted alphabetically
# by the "file" key so that file/r_count/currency/sum/built_in_total line up
# correctly again.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("87811004_1121_AU", 307, "AUD", 1617.46, "1617.46"),
    @("87811004_1121_BG", 6, "BGN", 16.04, "16.04"),
    @("87811004_1121_BR", 38, "BRL", 440.44, "440.44"),
    @("87811004_1121_CA", 287, "CAD", 1332.8, "1332.8"),
    @("87811004_1121_CH", 78, "CHF", 387.56, "387.56"),
    @("87811004_1121_CL", 35, "CLP", 73248, "73248"),
    @("87811004_1121_CO", 31, "COP", 347830, "347830"),
    @("87811004_1121_CZ", 23, "CZK", 1718.21, "1718.21"),
    @("87811004_1121_DK", 22, "DKK", 440.16, "440.16"),
    @("87811004_1121_EU", 746, "EUR", 3309.31, "3309.31"),
    @("87811004_1121_GB", 400, "GBP", 1309.43, "1309.43"),
    @("87811004_1121_HU", 900, "HUF", 1915878, "1915878"),
    @("87811004_1121_JP", 23, "JPY", 7546, "7546"),
    @("87811004_1121_LL", 43, "USD", 112, "112"),
    @("87811004_1121_MX", 78, "MXN", 7567.7, "7567.7"),
    @("87811004_1121_NO", 26, "NOK", 655.2, "655.2"),
    @("87811004_1121_NZ", 34, "NZD", 147.6, "147.6"),
    @("87811004_1121_PE", 21, "PEN", 174.3, "174.3"),
    @("87811004_1121_PL", 42, "PLN", 401.94, "401.94"),
    @("87811004_1121_RO", 494, "RON", 9795.57, "9795.57"),
    @("87811004_1121_SE", 36, "SEK", 919.87, "919.87"),
    @("87811004_1121_US", 1415, "USD", 8030.4, "8030.4")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]
    $eCell.Style = "Normal"
    $r++
}
